$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date value in A3 (from 2024-03-12 / 45363 to 2024-03-22 / 45373)
$ws.Range("A3").Value = 45373

# Update the active selection to A3 (was B5)
$ws.Range("A3").Select()
